$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Add a new "AFN" header label above the existing NFA table (row 10 was blank) ---
$ws1.Range("A10").Value = "AFN"

# --- Rework the DFA table that used to start at row 15 ---
$ws1.Range("A15").Value = "AFD"
$ws1.Range("B15").ClearContents()
$ws1.Range("C15").ClearContents()
$ws1.Range("D15").ClearContents()

$ws1.Range("A16").Value = "inicial"
$ws1.Range("B16").Value = "q0q1q2"
# C16 ("q0q1q2") and D16 ("q1q2") stay as-is

$ws1.Range("A17").ClearContents()
# B17 ("q1q2"), C17 ("q2"), D17 ("q1q2") stay as-is

$ws1.Range("B18").Value = "q2"
$ws1.Range("C18").Value = "q2"
$ws1.Range("D18").ClearContents()

$ws1.Range("A19").ClearContents()
$ws1.Range("B19").ClearContents()
$ws1.Range("C19").ClearContents()

# --- Add the new blank "Sheet2" right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("A2").Select()

# --- Restore Sheet1 as the active tab with its saved cursor position ---
$ws1.Activate()
$ws1.Range("B11").Select()
